# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 1.54
$ws.Range("H2").Value = 1.28
$ws.Range("I2").Value = 9.4
$ws.Range("K2").Value = 5.3

# Row 3
$ws.Range("F3").Value = 2.1
$ws.Range("J3").Value = 3.45
$ws.Range("K3").Value = 4

# Row 5
$ws.Range("F5").Value = 1.4
$ws.Range("G5").Value = 1.48
$ws.Range("K5").Value = 5.7

# Row 6
$ws.Range("G6").Value = 3.2
$ws.Range("J6").Value = 3.45
$ws.Range("K6").Value = 4
$ws.Range("Q6").Value = 1.77

# Row 7
$ws.Range("F7").Value = 1.73
$ws.Range("G7").Value = 2.04
$ws.Range("H7").Value = 3.45
$ws.Range("P7").Value = 2.14

# Row 8
$ws.Range("F8").Value = 2.26
$ws.Range("G8").Value = 2.96
$ws.Range("H8").Value = 2.68
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 5.7

# Row 9
$ws.Range("F9").Value = 1.61
$ws.Range("G9").Value = 1.68
$ws.Range("H9").Value = 5.2
$ws.Range("I9").Value = 6.6
$ws.Range("J9").Value = 4.4
$ws.Range("K9").Value = 4.9
$ws.Range("P9").Value = 2.46
